# Insert a new row into the daily price log for "Zapallo italiano" (Macroferia
# Regional de Talca). This pushes the previous rows 181-289 down to 182-290
# and populates the newly inserted row 181 with a new daily record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 181 (shifts 181:289 -> 182:290)
$ws.Rows.Item(181).Insert()

# Populate the new row 181 with the new record's data.
$ws.Range("A181").Value = 5
$ws.Range("B181").Value = "Macroferia Regional de Talca"
$ws.Range("C181").Value = "Maule"
$ws.Range("D181").Value = 44606
$ws.Range("E181").Value = 7
$ws.Range("F181").Value = 100112032
$ws.Range("G181").Value = "Zapallo italiano"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 500
$ws.Range("K181").Value = 4000
$ws.Range("L181").Value = 4000
$ws.Range("M181").Value = 4000
$ws.Range("N181").Value = "`$/caja 50 unidades"
$ws.Range("O181").Value = "Región del Maule"
$ws.Range("P181").Value = 80
$ws.Range("Q181").Value = 50
$ws.Range("R181").Value = "Hortaliza"

# Make sure the date cell keeps the same date number format used by the
# rest of column D.
$ws.Range("D181").NumberFormat = $ws.Range("D182").NumberFormat
